$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Value = "divesh"
$ws.Range("B29").Value = "x"
$ws.Range("C29").Value = "x"
$ws.Range("D29").Value = "must street"
$ws.Range("E29").Value = "nilesh"
$ws.Range("F29").Value = "y"
$ws.Range("G29").Value = "y"
$ws.Range("H29").Value = "rust street"

$ws.Range("I29").NumberFormat = "@"
$ws.Range("I29").Value = "10.12.2019"
$ws.Range("I29").Style = "Normal"

$ws.Range("J29").NumberFormat = "@"
$ws.Range("J29").Value = "12000"
$ws.Range("J29").Style = "Normal"

$ws.Range("K29").Value = "nilesh"
$ws.Range("L29").Value = "yamaha"
$ws.Range("M29").Value = "y2"
$ws.Range("N29").Value = "TN-02-C-1234"

$ws.Range("O29").NumberFormat = "@"
$ws.Range("O29").Value = "13"
$ws.Range("O29").Style = "Normal"

$ws.Range("P29").NumberFormat = "@"
$ws.Range("P29").Value = "2.1"
$ws.Range("P29").Style = "Normal"

$ws.Range("Q29").NumberFormat = "@"
$ws.Range("Q29").Value = "1000"
$ws.Range("Q29").Style = "Normal"

$ws.Range("R29").NumberFormat = "@"
$ws.Range("R29").Value = "2"
$ws.Range("R29").Style = "Normal"

$ws.Range("S29").Value = "r,c"

$ws.Range("T29").NumberFormat = "@"
$ws.Range("T29").Value = "123"
$ws.Range("T29").Style = "Normal"

$ws.Range("U29").Value = "e2fcc026-d886-11e9-abcb-107d1a2a80c2"
